$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap cell formatting between B and C for rows where the "winner"
#     highlight (style index 2: bold, no border) moves to the other column.
#     Row 3 keeps its original formatting (B3 plain, C3 highlighted) so it
#     is left untouched below.

# Row 2: highlight moves from C2 to B2
$ws.Range("C2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").ClearFormats() | Out-Null

# Row 4: highlight moves from B4 to C4
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").ClearFormats() | Out-Null

# Row 5: highlight moves from C5 to B5
$ws.Range("C5").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").ClearFormats() | Out-Null

# Row 6: highlight moves from C6 to B6
$ws.Range("C6").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("C6").ClearFormats() | Out-Null

$excel.CutCopyMode = $false

# --- Update the underlying score values
$ws.Range("B2").Value = 0.3448735177516937
$ws.Range("C2").Value = 0.3289481997489929

$ws.Range("B3").Value = 0.4984879726753482
$ws.Range("C3").Value = 0.5159719049806877

$ws.Range("B4").Value = 0.2227258788774616
$ws.Range("C4").Value = 0.2120054008814109

$ws.Range("B5").Value = 0.2230000048875809
$ws.Range("C5").Value = 0.1959999948740005

$ws.Range("B6").Value = -0.2570435404777527
$ws.Range("C6").Value = -0.278084397315979
